$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")

# --- Bump published version / date (values only change, rows stay put) ---
$ws1.Range("B3").Value = "0.2.0"
$ws1.Range("B8").Value = "2023-10-20T08:59:58+00:00"

# --- Insert a new "Jurisdiction" property row right after "Contact" (row 10) ---
# Grab the formatting of an existing body row *before* touching the sheet,
# then push row 11 ("Description" ...) and everything below it down by one.
$ws1.Range("A10:B10").Copy()
$ws1.Rows.Item(11).Insert(-4121)
$ws1.Range("A11:B11").PasteSpecial(-4122)

# Fill in the new row's real content.
$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = "iso:code:3166:FR"

Write-Output "ok"
